$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: true -> (cleared, row now has a blank value cell)
$ws.Range("B7").ClearContents()

# Date: 2023-10-31 -> 2025-11-18
# Assigning a date-shaped string directly via .Value causes Excel to
# auto-convert it into a date serial number (and allocate a new date
# number-format style). Route the literal text through a formula +
# copy/paste-values round trip instead, so it lands back as plain text
# using the exact same cell style as before.
$ws.Range("Z1").Formula = "=""2025-11-18"""
$ws.Range("Z1").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
